$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("A477")
$cell.Value = "Translated <contribution-count> <language> <span>sentence(s)</span>"

$len1 = "Translated <contribution-count> <language> <span>".Length
$len2 = "sentence(s)".Length
$len3 = "</span>".Length

$run2 = $cell.Characters($len1 + 1, $len2)
$run2.Font.Size = 12
$run2.Font.Color = 15238730
$run2.Font.Name = "Arial"

$run3 = $cell.Characters($len1 + $len2 + 1, $len3)
$run3.Font.Name = "Arial"
$run3.Font.Size = 12
$run3.Font.ThemeColor = 1
Write-Host "done"
